$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.996.07"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "2.587.69"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "2.597.87"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("D14").Value = "3.048.03"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "58.941.69"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").Value = "2.609.53"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.404"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0₃0725"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.06%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.45%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.26%  "
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.590"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "1.963.55"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.15%  "
